$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Target cluster = ECs)
$ws.Range("G2").Value = 0.09179766666666667
$ws.Range("H2").Value = 0.275393
$ws.Range("M2").Value = 0.006768333333333334
$ws.Range("N2").Value = 0.020305
$ws.Range("O2").Value = 0.00347174015482542
$ws.Range("P2").Value = 0.00347174015482542
$ws.Range("Q2").Value = 0.0006213172072222223
$ws.Range("R2").Value = 0.005591854865
$ws.Range("S2").Value = 0.00347174015482542
$ws.Range("T2").Value = 0.00347174015482542

# Update row 3 (Target cluster = FAPs)
$ws.Range("G3").Value = 0.09179766666666667
$ws.Range("H3").Value = 0.275393
$ws.Range("O3").Value = 0.9965282598451746
$ws.Range("P3").Value = 0.9965282598451747
$ws.Range("Q3").Value = 0.1783428850412222
$ws.Range("R3").Value = 1.605085965371
$ws.Range("S3").Value = 0.9965282598451746
$ws.Range("T3").Value = 0.9965282598451747

# Delete row 4 (Target cluster = MuSCs) entirely
$ws.Rows.Item(4).Delete()
